# Scheduled market-data refresh for Spriggan_Profits.xlsx.
# Updates cached currentAveragePrice / LevePrice / LeveProfit columns
# (H,I,J,K,L,M,N) for the leves whose Universalis price snapshot changed,
# across all eight Disciple of the Hand sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 21.666666
$ws.Range("I8").Value = 21.666666
$ws.Range("K8").Value = 64.99999800000001
$ws.Range("M8").Value = 74.00000199999999
$ws.Range("H15").Value = 1280.5172
$ws.Range("I15").Value = 1280.5172
$ws.Range("K15").Value = 3841.5516
$ws.Range("M15").Value = -3672.5516
$ws.Range("H28").Value = 1033.8125
$ws.Range("I28").Value = 570.5454999999999
$ws.Range("K28").Value = 570.5454999999999
$ws.Range("M28").Value = -85.54549999999995
$ws.Range("H43").Value = 5399.875
$ws.Range("I43").Value = 6300
$ws.Range("K43").Value = 6300
$ws.Range("M43").Value = -6231
$ws.Range("H45").Value = 700
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1308
$ws.Range("H64").Value = 35717516
$ws.Range("I64").Value = 50002960
$ws.Range("J64").Value = 3899.5
$ws.Range("K64").Value = 50002960
$ws.Range("L64").Value = 3899.5
$ws.Range("M64").Value = -50002712
$ws.Range("N64").Value = -4395.5
$ws.Range("H67").Value = 35717516
$ws.Range("I67").Value = 50002960
$ws.Range("J67").Value = 3899.5
$ws.Range("K67").Value = 50002960
$ws.Range("L67").Value = 3899.5
$ws.Range("M67").Value = -50002102
$ws.Range("N67").Value = -5615.5
$ws.Range("H109").Value = 47500
$ws.Range("J109").Value = 47500
$ws.Range("L109").Value = 47500
$ws.Range("N109").Value = -50274
$ws.Range("H112").Value = 93389
$ws.Range("I112").Value = 168333.17
$ws.Range("J112").Value = 65284.938
$ws.Range("K112").Value = 504999.51
$ws.Range("L112").Value = 195854.814
$ws.Range("M112").Value = -503891.51
$ws.Range("N112").Value = -198070.814
$ws.Range("H116").Value = 5000.85
$ws.Range("I116").Value = 5000.85
$ws.Range("K116").Value = 5000.85
$ws.Range("M116").Value = -1558.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 670296.8
$ws.Range("I2").Value = 818473.25
$ws.Range("K2").Value = 818473.25
$ws.Range("M2").Value = -818360.25
$ws.Range("H32").Value = 1163.6538
$ws.Range("I32").Value = 1163.6538
$ws.Range("K32").Value = 1163.6538
$ws.Range("M32").Value = -876.6538
$ws.Range("H74").Value = 27030900
$ws.Range("I74").Value = 31253358
$ws.Range("K74").Value = 31253358
$ws.Range("M74").Value = -31252484
$ws.Range("H77").Value = 27030900
$ws.Range("I77").Value = 31253358
$ws.Range("K77").Value = 156266790
$ws.Range("M77").Value = -156262422
$ws.Range("H96").Value = 39955.5
$ws.Range("J96").Value = 39955.5
$ws.Range("L96").Value = 39955.5
$ws.Range("N96").Value = -45447.5
$ws.Range("H116").Value = 670296.8
$ws.Range("I116").Value = 818473.25
$ws.Range("K116").Value = 818473.25
$ws.Range("M116").Value = -816179.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 670296.8
$ws.Range("I3").Value = 818473.25
$ws.Range("K3").Value = 818473.25
$ws.Range("M3").Value = -818359.25
$ws.Range("H64").Value = 750.5714
$ws.Range("J64").Value = 767.5
$ws.Range("L64").Value = 767.5
$ws.Range("N64").Value = -1217.5
$ws.Range("H67").Value = 750.5714
$ws.Range("J67").Value = 767.5
$ws.Range("L67").Value = 767.5
$ws.Range("N67").Value = -2327.5
$ws.Range("H99").Value = 2084.2856
$ws.Range("I99").Value = 2042
$ws.Range("K99").Value = 2042
$ws.Range("M99").Value = -544
$ws.Range("H105").Value = 2259.4
$ws.Range("I105").Value = 2299.25
$ws.Range("K105").Value = 2299.25
$ws.Range("M105").Value = -552.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5194.52
$ws.Range("I31").Value = 3139.5881
$ws.Range("J31").Value = 9561.25
$ws.Range("K31").Value = 3139.5881
$ws.Range("L31").Value = 9561.25
$ws.Range("M31").Value = -2844.5881
$ws.Range("N31").Value = -10151.25
$ws.Range("H34").Value = 5194.52
$ws.Range("I34").Value = 3139.5881
$ws.Range("J34").Value = 9561.25
$ws.Range("K34").Value = 3139.5881
$ws.Range("L34").Value = 9561.25
$ws.Range("M34").Value = -2937.5881
$ws.Range("N34").Value = -9965.25
$ws.Range("H35").Value = 1439
$ws.Range("I35").Value = 1148.3334
$ws.Range("J35").Value = 1875
$ws.Range("K35").Value = 1148.3334
$ws.Range("L35").Value = 1875
$ws.Range("M35").Value = -854.3334
$ws.Range("N35").Value = -2463
$ws.Range("H134").Value = 12502710
$ws.Range("I134").Value = 16669363
$ws.Range("K134").Value = 50008089
$ws.Range("M134").Value = -50005554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 3024.5
$ws.Range("I24").Value = 50
$ws.Range("J24").Value = 5999
$ws.Range("K24").Value = 150
$ws.Range("L24").Value = 17997
$ws.Range("M24").Value = 80
$ws.Range("N24").Value = -18457
$ws.Range("H94").Value = 18102.2
$ws.Range("I94").Value = 5498.75
$ws.Range("K94").Value = 16496.25
$ws.Range("M94").Value = -15820.25
$ws.Range("H131").Value = 1358.7826
$ws.Range("I131").Value = 979.7646999999999
$ws.Range("J131").Value = 2432.6667
$ws.Range("K131").Value = 2939.2941
$ws.Range("L131").Value = 7298.000100000001
$ws.Range("M131").Value = 2100.7059
$ws.Range("N131").Value = -17378.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6958.8
$ws.Range("I70").Value = 7112.7144
$ws.Range("J70").Value = 6599.6665
$ws.Range("K70").Value = 7112.7144
$ws.Range("L70").Value = 6599.6665
$ws.Range("M70").Value = -6842.7144
$ws.Range("N70").Value = -7139.6665
$ws.Range("H73").Value = 6958.8
$ws.Range("I73").Value = 7112.7144
$ws.Range("J73").Value = 6599.6665
$ws.Range("K73").Value = 7112.7144
$ws.Range("L73").Value = 6599.6665
$ws.Range("M73").Value = -6176.7144
$ws.Range("N73").Value = -8471.666499999999
$ws.Range("H80").Value = 3322.4443
$ws.Range("I80").Value = 3322.4443
$ws.Range("K80").Value = 3322.4443
$ws.Range("M80").Value = -2324.4443
$ws.Range("H83").Value = 3322.4443
$ws.Range("I83").Value = 3322.4443
$ws.Range("K83").Value = 16612.2215
$ws.Range("M83").Value = -11620.2215
$ws.Range("H107").Value = 2129.3333
$ws.Range("I107").Value = 568.25
$ws.Range("K107").Value = 568.25
$ws.Range("M107").Value = 1351.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1300
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1300
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 16635989
$ws.Range("I100").Value = 18148078
$ws.Range("K100").Value = 18148078
$ws.Range("M100").Value = -18147537
$ws.Range("H132").Value = 7580825.5
$ws.Range("I132").Value = 8934046
$ws.Range("K132").Value = 26802138
$ws.Range("M132").Value = -26799608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5298.3335
$ws.Range("J62").Value = 6777.222
$ws.Range("L62").Value = 6777.222
$ws.Range("N62").Value = -8025.222
$ws.Range("H65").Value = 5298.3335
$ws.Range("J65").Value = 6777.222
$ws.Range("L65").Value = 33886.11
$ws.Range("N65").Value = -40126.11
$ws.Range("H96").Value = 1549.875
$ws.Range("I96").Value = 1542.7142
$ws.Range("K96").Value = 1542.7142
$ws.Range("M96").Value = -169.7141999999999
$ws.Range("H136").Value = 11365482
$ws.Range("J136").Value = 1158.2858
$ws.Range("L136").Value = 3474.8574
$ws.Range("N136").Value = -8574.857400000001
